# GroupGrading.xlsx update
# - Grades sheet: add guide heading / group number / name rows, wrap header
#   labels onto two lines, and fill in the three group-member grading rows
#   (scores, totals, notes).
# - Grading Guide sheet: no content changes (shared-string indices shift
#   automatically as a consequence of the Grades-sheet string edits).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Grades")
$ws.Activate()

$nl = [char]10

# --- Apply wrap-text formatting first so it folds into the existing
#     header/body cell styles instead of resetting them. ---------------------
$ws.Range("B6:E6").WrapText = $true
$ws.Range("E7:E10").WrapText = $true

# --- Row 2/3: guide title + group number + student name -------------------
$ws.Range("A2").Value = "TEAM GRADING GUIDE"
$ws.Range("A3").Value = "Your Group Number: "
$ws.Range("B3").Value = "James Laurence - w0211593"

# --- Row 6: header labels, now wrapped onto two lines ----------------------
$ws.Range("B6").Value = "Effort and Commitment " + $nl + "[Teamwork] (10pts)"
$ws.Range("C6").Value = "Skills and " + $nl + "Knowedge (5pts)"
$ws.Range("D6").Value = "Project Management" + $nl + "and Coordination (5pts)"
$ws.Range("E6").Value = "Notes"
$ws.Range("F6").Value = "Total"

$ws.Rows.Item(6).RowHeight = 28.8

# --- Rows 7-9: group member grading data ------------------------------------
$ws.Range("A7").Value = "Gabriela Mkonde"
$ws.Range("B7").Value = 9
$ws.Range("C7").Value = 7
$ws.Range("D7").Value = 8
$ws.Range("E7").Value = "Gabby, was always reliable to complete tasks within the timeframe, or provide reason as to why she couldn't."
$ws.Range("F7").Formula = "=SUM(B7:D7)"
$ws.Rows.Item(7).RowHeight = 28.8

$ws.Range("A8").Value = "Louise Fear"
$ws.Range("B8").Value = 8
$ws.Range("C8").Value = 7
$ws.Range("D8").Value = 8
$ws.Range("E8").Value = "With Louise being sick for a few weeks, I was unable to rely on her for the tasks I delegated and had to increase my workload to compensate this."
$ws.Range("F8").Formula = "=SUM(B8:D8)"
$ws.Rows.Item(8).RowHeight = 43.2

$ws.Range("A9").Value = "Chris Whalen"
$ws.Range("B9").Value = 5
$ws.Range("C9").Value = 2
$ws.Range("D9").Value = 2
$ws.Range("E9").Value = "Chris did not show any initiative, had to wait to be told to do something. He seemed uninterested in doing anything outside his comfort zone of the material that was previously covered over the program. When I asked him to do something, he would try but would give up. (i.e. the xml extraction). With the minimum tasks that he did have, the time he spent on them was excessive, even though we were on a holding pattern due to the data extraction issue. he could have stepped up and offer to do other tasks, or even start something new. "
$ws.Range("F9").Formula = "=SUM(B9:D9)"
$ws.Rows.Item(9).RowHeight = 115.2

# --- Column widths -----------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 42.11
$ws.Columns.Item(2).ColumnWidth = 24.22
$ws.Columns.Item(3).ColumnWidth = 15.33
$ws.Columns.Item(4).ColumnWidth = 21.33
$ws.Columns.Item(5).ColumnWidth = 63.33

# --- Selection, matching the saved file's cursor position ------------------
$ws.Range("B9").Select()

Write-Output "Grades sheet updated"
